$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preconditions text reused from existing row (shared across the sheet)
$preconditions = "1-Go to (Travel-Adviser-Web-Application/UI/Log-Reg-UI/login.php)`n2- Enter valid Admin Credentials"

# ---------------------------------------------------------------------------
# Step 1: enter the three Bug IDs first (A9:A11)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "TAWA_Admin_Bug_008"
$ws.Range("A10").Value = "TAWA_Admin_Bug_009"
$ws.Range("A11").Value = "TAWA_Admin_Bug_010"

# ---------------------------------------------------------------------------
# Step 2: fill in row 9 (logout bug)
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "TAWA_Admin_019"
$ws.Range("D9").Value = $preconditions
$ws.Range("E9").Value = "Press logout"
$ws.Range("F9").Value = " The ‘admin shall be redirected to login page"
$ws.Range("G9").Value = "No action happens after pressing logout"
$ws.Range("C9").Value = "No action happens after pressing logout"
$ws.Range("H9").Value = "High"
$ws.Range("I9").Value = "High"
$ws.Range("J9").Value = "Functional"
$ws.Range("K9").Value = "Open"
$ws.Range("L9").Value = "Asmaa Hamdy"

$ws.Range("D9").WrapText = $true
$ws.Range("E9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 75

# ---------------------------------------------------------------------------
# Step 3: fill in row 10 (update username bug)
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "TAWA_Admin_020"
$ws.Range("D10").Value = $preconditions
$ws.Range("E10").Value = "1- Type an existing username in the search bar`n2- Press search`n3- Press Edit`n4- Update the current username to a new one`n5- Press update`n6- Press back to ""all users Page""`n7- Search for the new user name"
$ws.Range("F10").Value = "The new shall exist in the list"
$ws.Range("G10").Value = "The new user doesn't exist in the list"
$ws.Range("H10").Value = "High"
$ws.Range("I10").Value = "High"
$ws.Range("J10").Value = "Functional"
$ws.Range("K10").Value = "Open"
$ws.Range("L10").Value = "Asmaa Hamdy"

$ws.Range("D10").WrapText = $true
$ws.Range("E10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 135

# ---------------------------------------------------------------------------
# Step 4: fill in row 11 (update user, press update bug)
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "TAWA_Admin_021"
$ws.Range("D11").Value = $preconditions
$ws.Range("E11").Value = "1- Type an existing username in the search bar`n2- Press search`n3- Press Edit`n4- Update the current username to a new one`n5-Press Update"
$ws.Range("F11").Value = "The admin shall be redirected to ""successfully added"" page"
$ws.Range("G11").Value = "No action happens after pressing update"
$ws.Range("H11").Value = "Medium"
$ws.Range("I11").Value = "Medium"
$ws.Range("J11").Value = "Functional"
$ws.Range("K11").Value = "Open"
$ws.Range("L11").Value = "Asmaa Hamdy"

$ws.Range("D11").WrapText = $true
$ws.Range("E11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 105

# ---------------------------------------------------------------------------
# Step 5: go back and fill in the Description column (C) for rows 10 and 11
# ---------------------------------------------------------------------------
$ws.Range("C10").Value = "The new user doesn't exist in the list after editing"
$ws.Range("C11").Value = "No action happens after pressing update user"

# ---------------------------------------------------------------------------
# Update the view: scroll window and select A14 as in the authored workbook
# ---------------------------------------------------------------------------
$ws.Range("A14").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
